$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = "Dr. Shimaa Ahmad Mekki, Dr. Majorelle Magdy, Dr. Servinaz Sayed Mohammad, Dr. Nourhan Mahmoud, Dr. Rana Abo-Zaid"
$ws.Cells.Item(4, 7).Value = "Dr. Shimaa Ahmad Mekki, Dr. Alshimaa Atef, Dr. Amira Sobhy, Dr. Heba Mahmoud Ali, Dr. Hend Mahmoud, Dr. Menna tuâ€™Allah Medhat"
$ws.Cells.Item(5, 7).Value = "Dr. Abeer Ragab, Dr. Fatma Elhady, Dr. Menna tu'Alllah Mohammad, Dr. Nada Gouda"
$ws.Cells.Item(6, 7).Value = "Dr. Kerelos Zareef, Dr. Nada Mohammad"
$ws.Cells.Item(8, 7).Value = "Dr. Aya Saeed, Dr. Amal Awwad"
$ws.Cells.Item(10, 7).Value = "Dr. Amany Raafat, Dr. Madeha Saeed, Dr. Esraa Mostafa, Dr. Arwa Al-Sayed, Dr. Marina Youhanna, Dr. Maryam Ahmad"
$ws.Cells.Item(12, 7).Value = "Dr. Sarah Mahdy, Dr. Nouran Mahmoud"
$ws.Cells.Item(18, 7).Value = "Dr. Aya Hanafy, Dr. Remon, Dr. Shorok Mohammad, Dr. Yasmin"
$ws.Cells.Item(19, 7).Value = "Dr. Naema Gomaa, Dr. Yassmen Ahmad, Dr. Salma Hassan, Dr. Neveen Nashaat, Dr. Nardine, Dr. Remon, Dr. Monica, Dr. Maryam Ashraf, Dr. Wafaa Ebida"
$ws.Cells.Item(20, 7).Value = "Dr. Yassmen Ahmad, Dr. Youstina Magdy, Dr. Marina Sorial, Dr. Aya Hanafy, Dr. Nardine, Dr. Remon, Dr. Wafaa Ebida"
$ws.Cells.Item(21, 7).Value = "Dr. Yasmin, Dr. Yassmen Ahmad, Dr. Shorok Mohammad, Dr. Neveen Nashaat, Dr. Monica"
$ws.Cells.Item(22, 7).Value = "Dr. Naema Gomaa, Dr. Remon, Dr. Monica, Dr. Wafaa Ebida"
$ws.Cells.Item(24, 7).Value = "Dr. Shimaa Ahmad Mekki, Dr. Majorelle Magdy, Dr. Servinaz Sayed Mohammad, Dr. Nourhan Mahmoud, Dr. Rana Abo-Zaid"
$ws.Cells.Item(25, 7).Value = "Administrator, Dr. Alshimaa Atef, Dr. Manar Montaser, Dr. Gehan Adel"
$ws.Cells.Item(26, 7).Value = "Dr. Shimaa Ahmad Mekki, Dr. Alshimaa Atef, Dr. Amira Sobhy, Dr. Heba Mahmoud Ali, Dr. Hend Mahmoud, Dr. Menna tuâ€™Allah Medhat"
$ws.Cells.Item(27, 7).Value = "Dr. Abeer Ragab, Dr. Fatma Elhady, Dr. Menna tu'Alllah Mohammad, Dr. Nada Gouda"
$ws.Cells.Item(28, 7).Value = "Dr. Kerelos Zareef, Dr. Nada Mohammad"
$ws.Cells.Item(30, 7).Value = "Dr. Aya Saeed, Dr. Amal Awwad"
$ws.Cells.Item(32, 7).Value = "Dr. Amany Raafat, Dr. Madeha Saeed, Dr. Esraa Mostafa, Dr. Arwa Al-Sayed, Dr. Marina Youhanna, Dr. Maryam Ahmad"
$ws.Cells.Item(34, 7).Value = "Dr. Sarah Mahdy, Dr. Nouran Mahmoud"
$ws.Cells.Item(40, 7).Value = "Dr. Aya Hanafy, Dr. Remon, Dr. Shorok Mohammad, Dr. Yasmin"
$ws.Cells.Item(41, 7).Value = "Dr. Naema Gomaa, Dr. Yassmen Ahmad, Dr. Salma Hassan, Dr. Neveen Nashaat, Dr. Nardine, Dr. Remon, Dr. Monica, Dr. Maryam Ashraf, Dr. Wafaa Ebida"
$ws.Cells.Item(42, 7).Value = "Dr. Yassmen Ahmad, Dr. Youstina Magdy, Dr. Marina Sorial, Dr. Aya Hanafy, Dr. Nardine, Dr. Remon, Dr. Wafaa Ebida"
$ws.Cells.Item(43, 7).Value = "Dr. Yasmin, Dr. Yassmen Ahmad, Dr. Shorok Mohammad, Dr. Neveen Nashaat, Dr. Monica"
$ws.Cells.Item(44, 7).Value = "Dr. Naema Gomaa, Dr. Remon, Dr. Monica, Dr. Wafaa Ebida"
$ws.Cells.Item(46, 7).Value = "Dr. Shimaa Ahmad Mekki, Dr. Nahla Nagiub, Dr. Nourhan Mahmoud, Dr. Hend Mahmoud"
$ws.Cells.Item(48, 7).Value = "Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad, Dr. Nourhan Mahmoud, Dr. Nahla Nagiub, Dr. Menna tuâ€™Allah Medhat"
$ws.Cells.Item(49, 7).Value = "Dr. Menna tu'Alllah Mohammad, Dr. Nada Gouda, Dr. Amera Ahmad Saad"
$ws.Cells.Item(54, 7).Value = "Dr. Mai Mustafa, Dr. Amany Raafat, Dr. Madeha Saeed, Dr. Merna Said, Dr. Basma Hamed, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya, Dr. Yasmeena Fattoh, Dr. Marwa Mustafa, Dr. Marina Youhanna, Dr. Maryam Ahmad"
$ws.Cells.Item(58, 7).Value = "Dr. Afaf Abdallah, Dr. Amr Saeed"
$ws.Cells.Item(59, 7).Value = "Dr. Enas Omran, Dr. Marian Samir, Dr. Walaa Ghanima"
$ws.Cells.Item(60, 7).Value = "Dr. Nancy Abd Al-Shafy, Dr. Amr Saeed"
$ws.Cells.Item(62, 7).Value = "Dr. Yassmen Ahmad, Dr. Aya Hanafy, Dr. Shorok Mohammad, Dr. Wafaa Ebida"
$ws.Cells.Item(63, 7).Value = "Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Maryam Ashraf, Dr. Wafaa Ebida"
$ws.Cells.Item(64, 7).Value = "Dr. Neveen Nashaat, Dr. Youstina Magdy, Dr. Wafaa Ebida"
$ws.Cells.Item(65, 7).Value = "Dr. Remon, Dr. Salma Hassan, Dr. Shorok Mohammad, Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Nardine, Dr. Ola Abd Al-Fattah, Dr. Aya Hanafy"
$ws.Cells.Item(66, 7).Value = "Dr. Marina Sorial, Dr. Aya Hanafy, Dr. Monica, Dr. Eman Mohammad Al, Dr. Maryam Ashraf"
$ws.Cells.Item(68, 7).Value = "Dr. Shimaa Ahmad Mekki, Dr. Nahla Nagiub, Dr. Nourhan Mahmoud, Dr. Hend Mahmoud"
$ws.Cells.Item(70, 7).Value = "Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad, Dr. Nourhan Mahmoud, Dr. Nahla Nagiub, Dr. Menna tuâ€™Allah Medhat"
$ws.Cells.Item(71, 7).Value = "Dr. Menna tu'Alllah Mohammad, Dr. Nada Gouda, Dr. Amera Ahmad Saad"
$ws.Cells.Item(76, 7).Value = "Dr. Mai Mustafa, Dr. Amany Raafat, Dr. Madeha Saeed, Dr. Merna Said, Dr. Basma Hamed, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya, Dr. Yasmeena Fattoh, Dr. Marwa Mustafa, Dr. Marina Youhanna, Dr. Maryam Ahmad"
$ws.Cells.Item(80, 7).Value = "Dr. Afaf Abdallah, Dr. Amr Saeed"
$ws.Cells.Item(81, 7).Value = "Dr. Enas Omran, Dr. Marian Samir, Dr. Walaa Ghanima"
$ws.Cells.Item(82, 7).Value = "Dr. Nancy Abd Al-Shafy, Dr. Amr Saeed"
$ws.Cells.Item(84, 7).Value = "Dr. Yassmen Ahmad, Dr. Aya Hanafy, Dr. Shorok Mohammad, Dr. Wafaa Ebida"
$ws.Cells.Item(85, 7).Value = "Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Maryam Ashraf, Dr. Wafaa Ebida"
$ws.Cells.Item(86, 7).Value = "Dr. Neveen Nashaat, Dr. Youstina Magdy, Dr. Wafaa Ebida"
$ws.Cells.Item(87, 7).Value = "Dr. Remon, Dr. Salma Hassan, Dr. Shorok Mohammad, Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Nardine, Dr. Ola Abd Al-Fattah, Dr. Aya Hanafy"
$ws.Cells.Item(88, 7).Value = "Dr. Marina Sorial, Dr. Aya Hanafy, Dr. Monica, Dr. Eman Mohammad Al, Dr. Maryam Ashraf"
$ws.Cells.Item(90, 7).Value = "Dr. Shimaa Ahmad Mekki, Dr. Mohammad El-Tanany, Dr. Manar Montaser"
$ws.Cells.Item(92, 7).Value = "Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad, Dr. Nourhan Mahmoud, Dr. Nahla Nagiub, Dr. Menna tuâ€™Allah Medhat"
$ws.Cells.Item(93, 7).Value = "Dr. Menna tu'Alllah Mohammad, Dr. Fatma Elhady, Dr. Abeer Ragab, Dr. Amera Ahmad Saad"
$ws.Cells.Item(96, 7).Value = "Dr. Mariam Nour El-Din, Dr. Nourhan Mohammad, Dr. Sara Nabil, Dr. Amal Awwad"
$ws.Cells.Item(98, 7).Value = "Dr. Mai Mustafa, Dr. Amany Raafat, Dr. Madeha Saeed, Dr. Merna Said, Dr. Basma Hamed, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya, Dr. Yasmeena Fattoh, Dr. Marwa Mustafa, Dr. Marina Youhanna, Dr. Maryam Ahmad"
$ws.Cells.Item(104, 7).Value = "Dr. Nancy Abd Al-Shafy, Dr. Amr Saeed"
$ws.Cells.Item(106, 7).Value = "Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Nardine, Dr. Remon, Dr. Monica, Dr. Youstina Magdy"
$ws.Cells.Item(107, 7).Value = "Dr. Yassmen Ahmad, Dr. Neveen Nashaat, Dr. Aya Hanafy, Dr. Monica, Dr. Maryam Ashraf, Dr. Wafaa Ebida"
$ws.Cells.Item(108, 7).Value = "Dr. Yassmen Ahmad, Dr. Youstina Magdy, Dr. Marina Sorial, Dr. Aya Hanafy, Dr. Nardine, Dr. Remon, Dr. Wafaa Ebida"
$ws.Cells.Item(110, 7).Value = "Dr. Yassmen Ahmad, Dr. Monica, Dr. Wafaa Ebida"
$ws.Cells.Item(111, 7).Value = "Dr. Marina Atef, Dr. Naema Gomaa, Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Nourham Mostafa, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Wafaa Ebida"
$ws.Cells.Item(112, 7).Value = "Dr. Shimaa Ahmad Mekki, Dr. Mohammad El-Tanany, Dr. Manar Montaser"
$ws.Cells.Item(114, 7).Value = "Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad, Dr. Nourhan Mahmoud, Dr. Nahla Nagiub, Dr. Menna tuâ€™Allah Medhat"
$ws.Cells.Item(115, 7).Value = "Dr. Menna tu'Alllah Mohammad, Dr. Fatma Elhady, Dr. Abeer Ragab, Dr. Amera Ahmad Saad"
$ws.Cells.Item(118, 7).Value = "Dr. Mariam Nour El-Din, Dr. Nourhan Mohammad, Dr. Sara Nabil, Dr. Amal Awwad"
$ws.Cells.Item(120, 7).Value = "Dr. Mai Mustafa, Dr. Amany Raafat, Dr. Madeha Saeed, Dr. Merna Said, Dr. Basma Hamed, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya, Dr. Yasmeena Fattoh, Dr. Marwa Mustafa, Dr. Marina Youhanna, Dr. Maryam Ahmad"
$ws.Cells.Item(126, 7).Value = "Dr. Nancy Abd Al-Shafy, Dr. Amr Saeed"
$ws.Cells.Item(128, 7).Value = "Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Nardine, Dr. Remon, Dr. Monica, Dr. Youstina Magdy"
$ws.Cells.Item(129, 7).Value = "Dr. Yassmen Ahmad, Dr. Neveen Nashaat, Dr. Aya Hanafy, Dr. Monica, Dr. Maryam Ashraf, Dr. Wafaa Ebida"
$ws.Cells.Item(130, 7).Value = "Dr. Yassmen Ahmad, Dr. Youstina Magdy, Dr. Marina Sorial, Dr. Aya Hanafy, Dr. Nardine, Dr. Remon, Dr. Wafaa Ebida"
$ws.Cells.Item(131, 7).Value = "Dr. Nardine, Dr. Marina Atef"
$ws.Cells.Item(132, 7).Value = "Dr. Yassmen Ahmad, Dr. Monica, Dr. Wafaa Ebida"
$ws.Cells.Item(133, 7).Value = "Dr. Marina Atef, Dr. Naema Gomaa, Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Nourham Mostafa, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Wafaa Ebida"
$ws.Cells.Item(134, 7).Value = "Dr. Amira Sobhy, Dr. Majorelle Magdy, Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Asmaa Reda"
$ws.Cells.Item(137, 7).Value = "Dr. Abeer Ragab, Dr. Fatma Elhady, Dr. Menna tu'Alllah Mohammad, Dr. Nada Gouda"
$ws.Cells.Item(140, 7).Value = "Dr. Aya Saeed, Dr. Amal Awwad"
$ws.Cells.Item(142, 7).Value = "Dr. Merna Said, Dr. Basma Hamed, Dr. Esraa Mostafa, Dr. Yasmeena Fattoh, Dr. Marwa Mustafa"
$ws.Cells.Item(144, 7).Value = "Dr. Mona Ibrahim Hussein, Dr. Khadija Osama"
$ws.Cells.Item(147, 7).Value = "Dr. Nourham Mostafa, Dr. Nancy Abd Al-Shafy"
$ws.Cells.Item(150, 7).Value = "Dr. Naema Gomaa, Dr. Yassmen Ahmad, Dr. Salma Hassan, Dr. Neveen Nashaat, Dr. Nardine, Dr. Remon, Dr. Monica, Dr. Maryam Ashraf, Dr. Wafaa Ebida"
$ws.Cells.Item(151, 7).Value = "Dr. Marina Atef, Dr. Yassmen Ahmad, Dr. Monica, Dr. Wafaa Ebida"
$ws.Cells.Item(153, 7).Value = "Dr. Marina Sorial, Dr. Aya Hanafy, Dr. Monica, Dr. Eman Mohammad Al, Dr. Maryam Ashraf"
$ws.Cells.Item(154, 7).Value = "Dr. Naema Gomaa, Dr. Remon, Dr. Salma Hassan, Dr. Wafaa Ebida"
$ws.Cells.Item(155, 7).Value = "Dr. Marina Atef, Dr. Naema Gomaa, Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Nourham Mostafa, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Wafaa Ebida"
$ws.Cells.Item(156, 7).Value = "Dr. Alshimaa Atef, Dr. Manar Montaser, Dr. Majorelle Magdy, Dr. Mohammad El-Tanany, Dr. Menna tuâ€™Allah Medhat"
$ws.Cells.Item(159, 7).Value = "Dr. Abeer Ragab, Dr. Fatma Elhady, Dr. Menna tu'Alllah Mohammad, Dr. Nada Gouda"
$ws.Cells.Item(162, 7).Value = "Dr. Aya Saeed, Dr. Amal Awwad"
$ws.Cells.Item(164, 7).Value = "Dr. Merna Said, Dr. Basma Hamed, Dr. Esraa Mostafa, Dr. Yasmeena Fattoh, Dr. Marwa Mustafa"
$ws.Cells.Item(165, 7).Value = "Dr. Sarah Mahdy, Dr. Nouran Mahmoud"
$ws.Cells.Item(166, 7).Value = "Dr. Mona Ibrahim Hussein, Dr. Khadija Osama"
$ws.Cells.Item(169, 7).Value = "Dr. Nourham Mostafa, Dr. Nancy Abd Al-Shafy"
$ws.Cells.Item(172, 7).Value = "Dr. Naema Gomaa, Dr. Yassmen Ahmad, Dr. Salma Hassan, Dr. Neveen Nashaat, Dr. Nardine, Dr. Remon, Dr. Monica, Dr. Maryam Ashraf, Dr. Wafaa Ebida"
$ws.Cells.Item(173, 7).Value = "Dr. Marina Atef, Dr. Yassmen Ahmad, Dr. Monica, Dr. Wafaa Ebida"
$ws.Cells.Item(175, 7).Value = "Dr. Marina Sorial, Dr. Aya Hanafy, Dr. Monica, Dr. Eman Mohammad Al, Dr. Maryam Ashraf"
$ws.Cells.Item(176, 7).Value = "Dr. Naema Gomaa, Dr. Remon, Dr. Salma Hassan, Dr. Wafaa Ebida"
$ws.Cells.Item(177, 7).Value = "Dr. Marina Atef, Dr. Naema Gomaa, Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Nourham Mostafa, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Wafaa Ebida"
